$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.569.36"
$ws.Range("E2").Value = "  -3.84%  "
$ws.Range("D3").Value = "'3.468.75"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'581.24"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").Value = "'172.09"
$ws.Range("E6").Value = "  -5.80%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  -3.80%  "
$ws.Range("D9").Value = "'3.467.01"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").Value = "'0.131"
$ws.Range("D11").Value = "'6.84"
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("E12").Value = "  -5.20%  "
$ws.Range("D13").Value = "'4.068.25"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "'29.89"
$ws.Range("E15").Value = "  -7.45%  "
$ws.Range("D16").Value = "'65.695.87"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("E17").Value = "  -4.29%  "
$ws.Range("D18").Value = "'3.472.23"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "'5.94"
$ws.Range("E19").Value = "  -4.67%  "
$ws.Range("D20").Value = "'13.92"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "'366.69"
$ws.Range("E21").Value = "  -7.52%  "
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "'72.11"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "'0.0000121"
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").Value = "'9.74"
$ws.Range("E27").Value = "  -6.85%  "
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "'24.12"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  -6.51%  "
$ws.Range("D32").Value = "'1.98"
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -9.26%  "
$ws.Range("D35").Value = "'7.08"
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "'29.52"
$ws.Range("E37").Value = "  +11.73%  "
$ws.Range("D38").Value = "'159.11"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").Value = "'0.884"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("E40").Value = "  -5.69%  "
$ws.Range("D41").Value = "'2.803.84"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "'4.43"
$ws.Range("E42").Value = "  -6.27%  "
$ws.Range("E43").Value = "  -11.90%  "
$ws.Range("D44").Value = "'6.33"
$ws.Range("E44").Value = "  -6.07%  "
$ws.Range("E45").Value = "  -5.04%  "
$ws.Range("D46").Value = "'39.86"
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("D47").Value = "'24.18"
$ws.Range("E47").Value = "  -8.83%  "
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("D49").Value = "'305.10"
$ws.Range("E49").Value = "  -8.24%  "
$ws.Range("D50").Value = "'0.822"
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.22"
$ws.Range("E51").Value = "  -2.66%  "
